$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected FilesTab (Neo4j) query in B4: drop the "File Type" and "Breed"
# columns from the RETURN clause (ICDC Breed 1-14 script fix).
$filesQuery = "MATCH (f:file)-->(parent)`nWITH DISTINCT f, parent`nMATCH (f)-[*]->(c:case)<--(demo:demographic)`nWHERE demo.breed IN ['Irish Setter']`nOPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`nOPTIONAL MATCH (samp:sample)-->(c)`nWITH DISTINCT f, parent, c, demo, diag, s`nRETURN  coalesce(f.file_name, '') AS ``File Name``,`n           coalesce(labels(parent)[0], '') AS ``Association``,`n        coalesce(f.file_description, '') AS ``Description``,`n        coalesce(f.file_format, '') AS ``Format``,`n        coalesce(f.file_size, '') AS ``Size``,`n        coalesce(c.case_id, '') AS ``Case ID``,`n         coalesce(diag.disease_term,'') AS Diagnosis , `n        coalesce(s.clinical_study_designation,'') AS ``Study Code``"

$ws.Range("B4").Value = $filesQuery

# Row 4 shrinks after the shorter query text.
$ws.Rows.Item(4).RowHeight = 217.5

# Selection/view moves to B4.
[void]$ws.Range("B4").Select()
